# análise descritiva - acréscimo das metas
#
# Insert 5 new "meta" columns (meta, meta_avg, meta_std, meta_min, meta_max)
# right after the existing "taxa_sucesso" column (column F), pushing every
# column from the old "arrecadado_sucesso" onward five positions to the
# right. Excel automatically shifts all the existing data/formatting when a
# real column insert is performed, so we only need to insert the columns and
# then populate the five new ones with their header text, number format and
# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 5 blank columns at G:K (old G:K = arrecadado_sucesso.. shifts to L:P, etc.)
$ws.Range("G1:K1").EntireColumn.Insert()

# 2) Header row (row 1) text for the new columns, matching style of the
#    other header cells (EntireColumn.Insert already carried the s="4"
#    header style into the new cells).
$ws.Range("G1").Value = "meta"
$ws.Range("H1").Value = "meta_avg"
$ws.Range("I1").Value = "meta_std"
$ws.Range("J1").Value = "meta_min"
$ws.Range("K1").Value = "meta_max"

# 3) Number format for the new data cells: same "R$ #,##0.00" currency
#    format used by the neighbouring monetary columns (style index 3).
$ws.Range("G2:K3").NumberFormat = "R$ #,##0.00"

# 4) Data values.
#    Row 2 ("flex" / "apoia.se" group) has no observations, all zeros.
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0

#    Row 3 ("flex" / "catarse" group) gets the computed meta statistics.
$ws.Range("G3").Value = 15599716.7029188
$ws.Range("H3").Value = 11279.62162177787
$ws.Range("I3").Value = 16430.30708090436
$ws.Range("J3").Value = 12.04441558726698
$ws.Range("K3").Value = 198811.9434626772
